# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the f930b1d3-a2e3-4b1d-a58f-bf0f6c4fac66 file.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: mark the row as failed for both locales ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet: the row's own Status column shares the same string ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handback transform failed"
$zhcn.Range("K3").Value = "Handback file name: wx5blse5.5mt is different with handoff file name: f930b1d3-a2e3-4b1d-a58f-bf0f6c4fac66.fee25f9320641895fff627e0d6893ea3746e4121.zh-cn."

# --- de-de sheet: the row's own Status column shares the same string ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handback transform failed"
$dede.Range("K3").Value = "Handback file name: wx5blse5.5mt is different with handoff file name: f930b1d3-a2e3-4b1d-a58f-bf0f6c4fac66.fee25f9320641895fff627e0d6893ea3746e4121.de-de."
